$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.264.10"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.568.40"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "585.04"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").Value = "148.53"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +3.99%  "
$ws.Range("E9").Value = "  +4.12%  "
$ws.Range("D10").Value = "5.65"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "3.027.56"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").Value = "63.170.91"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("E16").Value = "  +5.27%  "
$ws.Range("D17").Value = "2.555.70"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "11.36"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "343.52"
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("E20").Value = "  +3.71%  "
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "66.68"
$ws.Range("E23").Value = "  +3.30%  "
$ws.Range("D24").Value = "2.692.79"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "8.20"
$ws.Range("E27").Value = "  +13.67%  "
$ws.Range("D28").Value = "8.55"
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").Value = "1.49"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +7.98%  "
$ws.Range("D32").Value = "0.0₃0826"
$ws.Range("D33").Value = "460.62"
$ws.Range("E33").Value = "  +13.59%  "
$ws.Range("D34").Value = "1.63"
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("D35").Value = "176.91"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").Value = "4.52"
$ws.Range("E38").Value = "  +4.90%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "151.48"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").Value = "21.14"
$ws.Range("E44").Value = "  +2.82%  "
$ws.Range("D45").Value = "0.0552"
$ws.Range("E45").Value = "  +7.01%  "
$ws.Range("D46").Value = "0.615"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").Value = "0.0981"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").Value = "0.0240"
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("D49").Value = "18.43"
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "11.38"
$ws.Range("E51").Value = "  -0.10%  "
